$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update Runmode column: B Suite -> N, E Suite -> Y, F Suite -> Y
$ws.Range("C3").Value = "N"
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"

# Update selection to C8
$ws.Range("C8").Select()
